# Scheduled-runner update: refresh market price / profit figures
# across the Titan_Profits leve-crafting worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Mirrors the authoritative per-cell OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 4490
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4490
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 13470
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -13806
$ws.Range("H100").Value = 33335694
$ws.Range("I100").Value = 55556556
$ws.Range("J100").Value = 4400
$ws.Range("K100").Value = 55556556
$ws.Range("L100").Value = 4400
$ws.Range("M100").Value = -55556015
$ws.Range("N100").Value = -5482
$ws.Range("H113").Value = 5577.8184
$ws.Range("I113").Value = 3078.5715
$ws.Range("K113").Value = 3078.5715
$ws.Range("M113").Value = 175.4285

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 21302
$ws.Range("I32").Value = 4199.9614
$ws.Range("K32").Value = 4199.9614
$ws.Range("M32").Value = -3912.9614
$ws.Range("H45").Value = 833.7
$ws.Range("I45").Value = 870.7778
$ws.Range("J45").Value = 500
$ws.Range("K45").Value = 870.7778
$ws.Range("L45").Value = 500
$ws.Range("M45").Value = -493.7778
$ws.Range("N45").Value = -1254
$ws.Range("H61").Value = 3529.8333
$ws.Range("I61").Value = 2439.5
$ws.Range("J61").Value = 5710.5
$ws.Range("K61").Value = 2439.5
$ws.Range("L61").Value = 5710.5
$ws.Range("M61").Value = -2227.5
$ws.Range("N61").Value = -6134.5
$ws.Range("H74").Value = 10438.714
$ws.Range("I74").Value = 1316.2858
$ws.Range("J74").Value = 19561.143
$ws.Range("K74").Value = 1316.2858
$ws.Range("L74").Value = 19561.143
$ws.Range("M74").Value = -442.2858000000001
$ws.Range("N74").Value = -21309.143
$ws.Range("H77").Value = 10438.714
$ws.Range("I77").Value = 1316.2858
$ws.Range("J77").Value = 19561.143
$ws.Range("K77").Value = 6581.429
$ws.Range("L77").Value = 97805.715
$ws.Range("M77").Value = -2213.429
$ws.Range("N77").Value = -106541.715
$ws.Range("H110").Value = 1385.5294
$ws.Range("I110").Value = 1237.8889
$ws.Range("J110").Value = 1551.625
$ws.Range("K110").Value = 1237.8889
$ws.Range("L110").Value = 1551.625
$ws.Range("M110").Value = 807.1111000000001
$ws.Range("N110").Value = -5641.625
$ws.Range("H122").Value = 8154.6
$ws.Range("I122").Value = 7909.2
$ws.Range("J122").Value = 8400
$ws.Range("K122").Value = 23727.6
$ws.Range("L122").Value = 25200
$ws.Range("M122").Value = -21277.6
$ws.Range("N122").Value = -30100
$ws.Range("H132").Value = 2711.5898
$ws.Range("I132").Value = 2239.3333
$ws.Range("J132").Value = 5309
$ws.Range("K132").Value = 6717.999899999999
$ws.Range("L132").Value = 15927
$ws.Range("M132").Value = -4187.999899999999
$ws.Range("N132").Value = -20987
$ws.Range("H136").Value = 3529.8333
$ws.Range("I136").Value = 2439.5
$ws.Range("J136").Value = 5710.5
$ws.Range("K136").Value = 7318.5
$ws.Range("L136").Value = 17131.5
$ws.Range("M136").Value = -4768.5
$ws.Range("N136").Value = -22231.5

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H105").Value = 373997.06
$ws.Range("I105").Value = 3675
$ws.Range("K105").Value = 3675
$ws.Range("M105").Value = -1928
$ws.Range("H134").Value = 4783.1875
$ws.Range("I134").Value = 3268.611
$ws.Range("J134").Value = 6730.5
$ws.Range("K134").Value = 9805.832999999999
$ws.Range("L134").Value = 20191.5
$ws.Range("M134").Value = -7270.832999999999
$ws.Range("N134").Value = -25261.5

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 5361.648
$ws.Range("I31").Value = 1842.4348
$ws.Range("J31").Value = 7972.6772
$ws.Range("K31").Value = 1842.4348
$ws.Range("L31").Value = 7972.6772
$ws.Range("M31").Value = -1547.4348
$ws.Range("N31").Value = -8562.6772
$ws.Range("H34").Value = 5361.648
$ws.Range("I34").Value = 1842.4348
$ws.Range("J34").Value = 7972.6772
$ws.Range("K34").Value = 1842.4348
$ws.Range("L34").Value = 7972.6772
$ws.Range("M34").Value = -1640.4348
$ws.Range("N34").Value = -8376.6772
$ws.Range("H58").Value = 2129.7026
$ws.Range("I58").Value = 1137.4445
$ws.Range("J58").Value = 4808.8
$ws.Range("K58").Value = 1137.4445
$ws.Range("L58").Value = 4808.8
$ws.Range("M58").Value = -934.4445000000001
$ws.Range("N58").Value = -5214.8
$ws.Range("H132").Value = 3262
$ws.Range("I132").Value = 1556
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 4668
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -2138
$ws.Range("N132").Value = -19964
$ws.Range("H134").Value = 2916.3635
$ws.Range("I134").Value = 1446.125
$ws.Range("J134").Value = 6837
$ws.Range("K134").Value = 4338.375
$ws.Range("L134").Value = 20511
$ws.Range("M134").Value = -1803.375
$ws.Range("N134").Value = -25581
$ws.Range("H136").Value = 2129.7026
$ws.Range("I136").Value = 1137.4445
$ws.Range("J136").Value = 4808.8
$ws.Range("K136").Value = 3412.3335
$ws.Range("L136").Value = 14426.4
$ws.Range("M136").Value = -862.3335000000002
$ws.Range("N136").Value = -19526.4
$ws.Range("H138").Value = 67050
$ws.Range("J138").Value = 67050
$ws.Range("L138").Value = 67050
$ws.Range("N138").Value = -77330

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H131").Value = 7577318
$ws.Range("I131").Value = 573.3333
$ws.Range("K131").Value = 1719.9999
$ws.Range("M131").Value = 3320.0001
$ws.Range("H132").Value = 1158.25
$ws.Range("I132").Value = 893.7857
$ws.Range("J132").Value = 1528.5
$ws.Range("K132").Value = 8044.071300000001
$ws.Range("L132").Value = 13756.5
$ws.Range("M132").Value = -5514.071300000001
$ws.Range("N132").Value = -18816.5
$ws.Range("H138").Value = 1349
$ws.Range("I138").Value = 822.7143
$ws.Range("J138").Value = 5033
$ws.Range("K138").Value = 2468.1429
$ws.Range("L138").Value = 15099
$ws.Range("M138").Value = 2671.8571
$ws.Range("N138").Value = -25379

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H132").Value = 3168.05
$ws.Range("I132").Value = 2960.1875
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 8880.5625
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -6350.5625
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H93").Value = 2829.125
$ws.Range("I93").Value = 3000
$ws.Range("J93").Value = 2726.6
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 2726.6
$ws.Range("M93").Value = -1752
$ws.Range("N93").Value = -5222.6
$ws.Range("H134").Value = 69429
$ws.Range("J134").Value = 69429
$ws.Range("L134").Value = 69429
$ws.Range("N134").Value = -79569

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 3527.1614
$ws.Range("I132").Value = 3456.3076
$ws.Range("J132").Value = 3895.6
$ws.Range("K132").Value = 10368.9228
$ws.Range("L132").Value = 11686.8
$ws.Range("M132").Value = -7838.9228
$ws.Range("N132").Value = -16746.8
